$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 54, shifting existing rows 54-91 down to 55-92.
$ws.Rows(54).Insert()

# Populate the newly inserted row 54 with the new weekly price record.
$ws.Cells.Item(54, 1).Value = 1
$ws.Cells.Item(54, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(54, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(54, 4).Value = 44574
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = 15
$ws.Cells.Item(54, 6).Value = 100112008
$ws.Cells.Item(54, 7).Value = "Coliflor"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Tercera"
$ws.Cells.Item(54, 10).Value = 1200
$ws.Cells.Item(54, 11).Value = 450
$ws.Cells.Item(54, 12).Value = 500
$ws.Cells.Item(54, 13).Value = 475
$ws.Cells.Item(54, 14).Value = "$/unidad"
$ws.Cells.Item(54, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(54, 16).Value = 475
$ws.Cells.Item(54, 17).Value = 1
$ws.Cells.Item(54, 18).Value = "Hortaliza"
